# ZBP_03_strategie_domacnosti.xlsx update
# - adds a new wave "22. 6. 2021" as the latest column on both sheets
# - refreshes the "aktualizace" date embedded in the two title rows

$wb = $excel.ActiveWorkbook

$wsData   = $wb.Worksheets.Item("data")
$wsPocetR = $wb.Worksheets.Item("pocetR")

# ---------------------------------------------------------------------
# Sheet "data": new column AE (31) = "22. 6. 2021", % values for rows 2-45
# ---------------------------------------------------------------------

# Give the new header cell the same formatting as the previous header (AD1)
$wsData.Cells.Item(1, 30).Copy($wsData.Cells.Item(1, 31))
$wsData.Cells.Item(1, 31).Value = "22. 6. 2021"

$dataValues = @(
    0.2, 0.12, 0.46, 0.29, 0.13, 0.19, 0.26, 0.17, 0.24,
    0.18, 0.23, 0.31, 0.19, 0.18, 0.23, 0.17, 0.26,
    0.25, 0.14, 0.13, 0.12, 0.23, 0.4, 0.43, 0.17,
    0.06, 0.12, 0.2, 0.08, 0.12, 0.12, 0.19, 0.19,
    0.09, 0.11, 0.15, 0.08, 0.25, 0.14, 0.07,
    0.07, 0.08, 0.18, 0.24
)

for ($i = 0; $i -lt $dataValues.Count; $i++) {
    $row = $i + 2
    $wsData.Cells.Item($row, 31).Value = $dataValues[$i]
}

# Title cell (row 46, column A) - bump the "aktualizace" date
$titleCell = $wsData.Cells.Item(46, 1)
$titleCell.Value = ($titleCell.Value2 -replace "1\. 6\. 2021", "28. 6. 2021")

# ---------------------------------------------------------------------
# Sheet "pocetR": new column AD (30) = "22. 6. 2021", counts for rows 2-23
# ---------------------------------------------------------------------

$wsPocetR.Cells.Item(1, 29).Copy($wsPocetR.Cells.Item(1, 30))
$wsPocetR.Cells.Item(1, 30).Value = "22. 6. 2021"

$pocetValues = @(
    1904, 183, 377, 1344, 911, 169, 545, 279, 866,
    159, 117, 762, 867, 661, 376, 194, 753,
    580, 249, 556, 345, 158
)

for ($i = 0; $i -lt $pocetValues.Count; $i++) {
    $row = $i + 2
    $wsPocetR.Cells.Item($row, 30).Value = $pocetValues[$i]
}

# Row 24 is the trailing title row; column AD needs the same "blank"
# cell shape as the rest of that row, so replicate it from AC24.
$wsPocetR.Cells.Item(24, 29).Copy($wsPocetR.Cells.Item(24, 30))

# Title cell (row 24, column A) - bump the "aktualizace" date
$titleCell2 = $wsPocetR.Cells.Item(24, 1)
$titleCell2.Value = ($titleCell2.Value2 -replace "1\. 6\. 2021", "28. 6. 2021")
